# Auto-generated edits updating the crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.007.06'
$ws.Range("E2").Value = '  +6.62%  '

$ws.Range("D3").Value = '3.010.72'
$ws.Range("E3").Value = '  +3.96%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.98'
$ws.Range("E5").Value = '  +2.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.14'
$ws.Range("E6").Value = '  +13.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").Value = '3.009.81'
$ws.Range("E8").Value = '  +3.97%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +3.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.00'
$ws.Range("E10").Value = '  +1.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  +7.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  +6.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +8.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.94'
$ws.Range("E14").Value = '  +8.63%  '

$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = '65.974.59'
$ws.Range("E16").Value = '  +6.73%  '

$ws.Range("D17").Value = '3.509.42'
$ws.Range("E17").Value = '  +3.97%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.96'
$ws.Range("E18").Value = '  +6.74%  '

$ws.Range("D19").Value = '3.009.45'
$ws.Range("E19").Value = '  +4.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '458.61'
$ws.Range("E20").Value = '  +6.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("E21").Value = '  +7.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.686'
$ws.Range("E22").Value = '  +4.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.37'
$ws.Range("E23").Value = '  +7.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.44'
$ws.Range("E24").Value = '  +4.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  +13.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.41'
$ws.Range("E26").Value = '  +2.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").Value = '  +5.66%  '

$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.14'
$ws.Range("E29").Value = '  +16.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.35'
$ws.Range("E30").Value = '  +15.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("E31").Value = '  +4.47%  '

$ws.Range("E32").Value = '  -6.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.99'
$ws.Range("E33").Value = '  +5.55%  '

$ws.Range("E34").Value = '  +3.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("E36").Value = '  +4.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.80'
$ws.Range("E37").Value = '  +7.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.15'
$ws.Range("E38").Value = '  +12.10%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  +6.37%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.92'
$ws.Range("E40").Value = '  +2.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.306'
$ws.Range("E41").Value = '  +14.15%  '

$ws.Range("E42").Value = '  +6.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.79'
$ws.Range("E43").Value = '  +8.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.49'
$ws.Range("E44").Value = '  +4.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '385.38'
$ws.Range("E45").Value = '  +11.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0356'
$ws.Range("E46").Value = '  +6.30%  '

$ws.Range("D47").Value = '2.796.54'
$ws.Range("E47").Value = '  +3.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.01'
$ws.Range("E48").Value = '  +2.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.00'
$ws.Range("E50").Value = '  +10.99%  '

$ws.Range("E51").Value = '  +4.17%  '

